# Update mock candidate data in place (rows 2-11, columns A-F)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("April Duncan", "william36@example.com", "Backend Developer", "SENIOR", "REJECTED", "Gilbertburgh"),
    @("Kyle Wheeler", "hhouse@example.com", "Backend Developer", "SENIOR", "OFFERED", "Ryanstad"),
    @("Mitchell Thompson", "dfaulkner@example.org", "Backend Developer", "SENIOR", "HIRED", "New Jackieview"),
    @("Robert Macias", "hernandezjoshua@example.org", "Full Stack Developer", "JUNIOR", "OFFERED", "Jenniferfort"),
    @("Cynthia Gonzales", "bguerrero@example.net", "UI/UX Designer", "SENIOR", "OFFERED", "Wellschester"),
    @("Richard Howard", "jenniferfreeman@example.org", "DevOps Engineer", "LEAD", "ON_HOLD", "Michaelbury"),
    @("Kelly Rose", "wschultz@example.com", "Backend Developer", "JUNIOR", "OFFER_ACCEPTED", "Jenniferland"),
    @("Alan Simmons", "dean13@example.net", "Product Manager", "LEAD", "SHORTLISTED", "Watkinshaven"),
    @("Kimberly Gonzales", "hparker@example.net", "Full Stack Developer", "JUNIOR", "REJECTED", "Nicoleville"),
    @("Ryan Dean", "steven07@example.org", "UI/UX Designer", "MID", "OFFERED", "Keithshire")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
}
